# Actualización automática de grupos experimentales
# Reassign the Grupo_Experimental (column B) values for rows 2-12,
# and populate row 12 (previously missing its group + had text "0.000"
# placeholders for SmartScore columns) with real values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Grupo_Experimental values per participant row.
$grupos = @{
    2  = "Sin SmartScore"
    3  = "Con SmartScore"
    4  = "Con SmartScore"
    5  = "Con SmartScore"
    6  = "Sin SmartScore"
    7  = "Con SmartScore"
    8  = "Sin SmartScore"
    10 = "Sin SmartScore"
    11 = "Con SmartScore"
    12 = "Sin SmartScore"
}

foreach ($row in $grupos.Keys) {
    $ws.Cells.Item($row, 2).Value = $grupos[$row]
}

# Row 12 (Paula Belén Chairez Rosas): the SmartScore columns were stored as
# text "0.000" placeholders; convert them to real numeric zeros.
$smartScoreCols = @(9, 12, 15, 18, 21, 24, 27, 30, 33)  # I, L, O, R, U, X, AA, AD, AG
foreach ($col in $smartScoreCols) {
    $ws.Cells.Item(12, $col).Value = 0
}
